# Remove the two paragraphs that were dropped from the assessment document:
#   1. "Specifics … exams? Quizzes? Question bank? "
#   2. "Colleen"
# Both paragraphs (including their paragraph marks) are deleted entirely.

$d = $word.ActiveDocument

function Remove-ParagraphByText($searchText) {
    $range = $d.Content
    $found = $range.Find.Execute($searchText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $para = $range.Paragraphs(1)
        $paraRange = $para.Range
        $paraRange.Delete()
    }
}

Remove-ParagraphByText("Specifics")
Remove-ParagraphByText("Colleen")
